# Append 5 new daily rows (234-238) to Sheet1, continuing the existing
# date series in column A (style copied from the last existing row so the
# date number format / border / alignment match) and zero values in B:D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 233
$newDates = @(44308, 44309, 44310, 44311, 44312)

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $r = $lastRow + 1 + $i

    # Copy the formatting of the template row (last existing data row) into
    # the new row so the new cells inherit identical styles (e.g. the date
    # cell style used in column A).
    $ws.Range("A$lastRow`:D$lastRow").Copy($ws.Range("A$r`:D$r"))

    $ws.Range("A$r").Value = $newDates[$i]
    $ws.Range("B$r").Value = 0
    $ws.Range("C$r").Value = 0
    $ws.Range("D$r").Value = 0
}

Write-Output "Added rows 234-238"
